$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "zobrazují souhvězdí Souhvězdí",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "zobrazujíSouhvězdí",
    2
)
